$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.01792803175822741, 0.6873742451724051, 1.225852550159523, 1.107182256974669, 1.137372497882221, 19),
    @(-0.06608967275348025, 0.7088863950418706, 1.217845040002488, 1.103560166009306, 1.133515860779062, 18),
    @(0.001870457576035796, 0.6187918867460556, 0.6230069128923086, 0.7893078695238688, 0.8135976448302396, 17),
    @(0.1204196123905843, 0.5277616828532142, 0.6040816134619066, 0.7772268738675385, 0.7930233996807294, 16),
    @(0.09828742272023362, 0.5953749172355302, 0.77910483569946, 0.8826691541565617, 0.9079673638532668, 15),
    @(0.1268913055883079, 0.6888184636477384, 1.036376752836633, 1.018025909707917, 1.048216613390359, 14),
    @(0.1309701764864839, 0.6044880343376112, 0.7446009131626561, 0.8629026093150119, 0.8877321498458313, 13),
    @(0.1991520669947553, 0.6396550887298748, 0.7959921877658299, 0.8921839427863684, 0.9083435325175242, 12),
    @(0.1894669389442902, 0.7879769448546992, 1.119193970586897, 1.057919642783372, 1.09161617549273, 11),
    @(0.1730349822307496, 0.6656591560163745, 0.9145633670025062, 0.9563280645272867, 0.9914200040160747, 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
